# Vega Modelo de Temuco - Alcachofa: weekly fruit/vegetable price update.
# Insert two new rows (119-120) above the current row 119, pushing the
# existing rows 119-129 down to 121-131, and populate the new rows with
# the latest week's price entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("119:120").Insert()

# Row 119: Vega Modelo de Temuco - Alcachofa - Española - Región Metropolitana
$ws.Cells.Item(119, 1).Value = 10
$ws.Cells.Item(119, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(119, 3).Value = "La Araucanía"
$ws.Cells.Item(119, 4).Value = 44474
$ws.Cells.Item(119, 5).Value = 9
$ws.Cells.Item(119, 6).Value = 100112013
$ws.Cells.Item(119, 7).Value = "Alcachofa"
$ws.Cells.Item(119, 8).Value = "Española"
$ws.Cells.Item(119, 9).Value = "Primera"
$ws.Cells.Item(119, 10).Value = 50
$ws.Cells.Item(119, 11).Value = 13000
$ws.Cells.Item(119, 12).Value = 13000
$ws.Cells.Item(119, 13).Value = 13000
$ws.Cells.Item(119, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(119, 15).Value = "Región Metropolitana"
$ws.Cells.Item(119, 16).Value = 433
$ws.Cells.Item(119, 17).Value = 30
$ws.Cells.Item(119, 18).Value = "Hortaliza"

# Row 120: Vega Modelo de Temuco - Alcachofa - Madrigal - Región Metropolitana
$ws.Cells.Item(120, 1).Value = 10
$ws.Cells.Item(120, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(120, 3).Value = "La Araucanía"
$ws.Cells.Item(120, 4).Value = 44474
$ws.Cells.Item(120, 5).Value = 9
$ws.Cells.Item(120, 6).Value = 100112013
$ws.Cells.Item(120, 7).Value = "Alcachofa"
$ws.Cells.Item(120, 8).Value = "Madrigal"
$ws.Cells.Item(120, 9).Value = "Primera"
$ws.Cells.Item(120, 10).Value = 50
$ws.Cells.Item(120, 11).Value = 12000
$ws.Cells.Item(120, 12).Value = 12000
$ws.Cells.Item(120, 13).Value = 12000
$ws.Cells.Item(120, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(120, 15).Value = "Región Metropolitana"
$ws.Cells.Item(120, 16).Value = 300
$ws.Cells.Item(120, 17).Value = 40
$ws.Cells.Item(120, 18).Value = "Hortaliza"
